$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Database sheet: insert 5 new rows (87-91) for new "(SK)" resources
# ---------------------------------------------------------------------------
$dbws = $wb.Worksheets.Item("Database")
$dbws.Rows("87:91").Insert()

# Row 87 - WasteWater (SK)
$dbws.Range("A87").Value = "WasteWater (SK)"
$dbws.Range("C87").Value = 1
$dbws.Range("E87").Formula = "=0.000625*1000"

# Row 88 - Water (SK)
$dbws.Range("A88").Value = "Water (SK)"
$dbws.Range("B88").Value = "H2O"
$dbws.Range("C88").Value = 1
$dbws.Range("D88").Value = 18.02
$dbws.Range("E88").Formula = "=0.000625*1000"
$dbws.Range("F88").Formula = "=C88*D88/(E88 * 1000)"

# Row 89 - Food (SK)
$dbws.Range("A89").Value = "Food (SK)"
$dbws.Range("E89").Formula = "= 0.000325*1000"

# Row 90 - Waste (SK)
$dbws.Range("A90").Value = "Waste (SK)"
$dbws.Range("C90").Value = 1
$dbws.Range("E90").Formula = "=0.000325*1000"

# Row 91 - Oxygen (SK)
$dbws.Range("A91").Value = "Oxygen (SK)"
$dbws.Range("B91").Value = "O2"
$dbws.Range("C91").Value = 1
$dbws.Range("D91").Value = 32
$dbws.Range("E91").Formula = "=0.0000013889*1000"
$dbws.Range("F91").Formula = "=C91*D91/(E91 * 1000)"

# ---------------------------------------------------------------------------
# 2. Fix the data-validation list ranges on sheet1/sheet2 (not auto-adjusted)
# ---------------------------------------------------------------------------
$kgws = $wb.Worksheets.Item("Calc (Kg)")
$kgws.Range("B7:F7").Validation.Delete()
$kgws.Range("B7:F7").Validation.Add(3, 1, 1, "=Database!`$A`$2:`$A`$128")
$kgws.Range("B14:F14").Validation.Delete()
$kgws.Range("B14:F14").Validation.Add(3, 1, 1, "=Database!`$A`$2:`$A`$128")

$molws = $wb.Worksheets.Item("Calc (Moles)")
$molws.Range("B16:F16").Validation.Delete()
$molws.Range("B16:F16").Validation.Add(3, 1, 1, "=Database!`$A`$2:`$A`$128")
$molws.Range("B7:F7").Validation.Delete()
$molws.Range("B7:F7").Validation.Add(3, 1, 1, "=Database!`$A`$2:`$A`$128")

# ---------------------------------------------------------------------------
# 3. Calc (Kg) sheet inputs
# ---------------------------------------------------------------------------
$kgws.Range("B6").Value = 0.000046296289999999999
$kgws.Range("B7").Value = "Waste (SK)"

$kgws.Range("B13").Value = 0.000024540774999999999
$kgws.Range("B14").Value = "Ammonia"
$kgws.Range("C14").Value = "Aerozine50"
$kgws.Range("D14").Value = "Water (SK)"
